$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting existing rows 84:222 down to 85:223
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record
$ws.Range("A84").Value = 3
$ws.Range("B84").Value = "Femacal de La Calera"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44540
$ws.Range("E84").Value = 5
$ws.Range("F84").Value = 100112039
$ws.Range("G84").Value = "Ciboulette"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 180
$ws.Range("K84").Value = 1500
$ws.Range("L84").Value = 1500
$ws.Range("M84").Value = 1500
$ws.Range("N84").Value = "$/docena de atados"
$ws.Range("O84").Value = "Provincia de Quillota"
$ws.Range("P84").Value = 500
$ws.Range("Q84").Value = 3
$ws.Range("R84").Value = "Hortaliza"
